$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 and Row 4 effectively swap their per-origin data (date, volume,
# min/max/avg price, origin name, price per kg) while the remaining
# descriptive columns stay the same.

$ws.Range("D2").Value = 44273
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 14000
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 233

$ws.Range("D4").Value = 44350
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = 10000
$ws.Range("O4").Value = "Región de Arica y Parinacota"
$ws.Range("P4").Value = 167
